$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.039.70"
$ws.Range("E2").Value = "  +3.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.895.81"
$ws.Range("E3").Value = "  +3.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.31"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.94"
$ws.Range("E8").Value = "  +2.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2942"
$ws.Range("E9").Value = "  +5.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06624"
$ws.Range("E10").Value = "  +3.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.899.20"
$ws.Range("E11").Value = "  +4.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.02"
$ws.Range("E12").Value = "  +1.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07231"
$ws.Range("E13").Value = "  +2.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6775"
$ws.Range("E14").Value = "  +5.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.24"
$ws.Range("E15").Value = "  +2.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.859"
$ws.Range("E16").Value = "  +3.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.026.48"
$ws.Range("E17").Value = "  +3.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007918"
$ws.Range("E18").Value = "  +8.11%  "

$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("E20").Value = "  +5.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.145.05"
$ws.Range("E21").Value = "  +5.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.770"
$ws.Range("E23").Value = "  +4.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.660"
$ws.Range("E24").Value = "  +5.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.211"
$ws.Range("E25").Value = "  +4.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.29"
$ws.Range("E26").Value = "  +1.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "131.44"
$ws.Range("E27").Value = "  +2.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.77"
$ws.Range("E28").Value = "  +2.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.964"
$ws.Range("E29").Value = "  +4.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.383"
$ws.Range("E30").Value = "  -1.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.223"
$ws.Range("E31").Value = "  +2.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08726"
$ws.Range("E32").Value = "  +4.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.940"
$ws.Range("E33").Value = "  +3.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05097"
$ws.Range("E34").Value = "  +3.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.122"
$ws.Range("E35").Value = "  +2.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7016"
$ws.Range("E36").Value = "  +4.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.679"
$ws.Range("E37").Value = "  -0.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.775"
$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.226"
$ws.Range("E39").Value = "  -3.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9524"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01658"
$ws.Range("E41").Value = "  +4.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.989"
$ws.Range("E42").Value = "  -2.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4219"
$ws.Range("E44").Value = "  +3.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.13"
$ws.Range("E45").Value = "  +2.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.470"
$ws.Range("E46").Value = "  +3.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1260"
$ws.Range("E47").Value = "  +2.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05745"
$ws.Range("E48").Value = "  +4.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.83"
$ws.Range("E49").Value = "  +3.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.221"
$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3730"
$ws.Range("E51").Value = "  +3.32%  "
